# "Se realizan cambios para sanity semilla 4"
# Update the seed/test data in the "Semilla 4" sheet (B9:B13) and flip which
# sheet/cell is the active selection when the workbook is reopened.

$wb = $excel.ActiveWorkbook

$wsSemilla4 = $wb.Worksheets.Item("Semilla 4")
$wsSemilla3 = $wb.Worksheets.Item("Semilla 3")

# Refresh the sanity seed values on "Semilla 4" (column B, rows 9-13).
# Written in this order so the new shared-string table entries land in the
# same sequence as the target workbook.
$wsSemilla4.Range("B10").Value = "874513783"
$wsSemilla4.Range("B11").Value = "680974538"
$wsSemilla4.Range("B12").Value = "1016141153"
$wsSemilla4.Range("B13").Value = "499757452"
$wsSemilla4.Range("B9").Value = "624917673"

# "Semilla 3" was previously the active/selected tab; make "Semilla 4" the
# active tab instead, with its selection moved to B18.
$wsSemilla4.Activate()
$wsSemilla4.Range("B18").Select()
